# Update the "Förändrad" (Changed) date column C for rows 2-28 from
# 45515 (2024-08-11) to 45516 (2024-08-12), keeping existing date formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45515) {
        $cell.Value2 = 45516
    }
}
